$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.931.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.380.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.959.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.384.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.070.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000114"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.521.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +7.28%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.416.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0761"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.778"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.429.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0260"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.83%  "
